$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (Oregon) updates ---
$ws.Range("B9").Value = 44237
$ws.Range("C9").Value = 148475
$ws.Range("D9").Value = 2044
$ws.Range("E9").Value = 3586
$ws.Range("H9").Value = 2.55
$ws.Range("K9").Value = 121163
$ws.Range("L9").Value = 1723

# --- Row 43 (Idaho) updates ---
$ws.Range("B43").NumberFormat = $ws.Range("B9").NumberFormat
$ws.Range("B43").Value = 44237
$ws.Range("C43").Value = 112993
$ws.Range("D43").Value = 1791
$ws.Range("E43").Value = 973
$ws.Range("F43").Value = 9
$ws.Range("G43").Value = 0.86
$ws.Range("H43").Value = 0.5
$ws.Range("J43").Value = $true
$ws.Range("O43").Value = "Success!"
